# Updated: po 31. 01. 2022
# Applies revised AgTests/AgPosit backfill values for late-Jan rows
# and appends three new daily rows (2022-01-28..2022-01-30).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column F only updates (rows with a revised AgTests figure)
$fOnlyUpdates = @{
    645 = 35739
    649 = 62665
    652 = 35206
    656 = 52558
    659 = 26386
    663 = 37257
    666 = 23912
    670 = 52624
    672 = 29825
    674 = 28712
    675 = 13487
    676 = 28092
    677 = 56184
    678 = 33823
    679 = 29360
}
foreach ($row in $fOnlyUpdates.Keys) {
    $ws.Cells.Item($row, 6).Value = $fOnlyUpdates[$row]
}

# Column F and G updates (rows with revised AgTests AND AgPosit figures)
$fgUpdates = @{
    680 = @(28371, 550)
    681 = @(26368, 576)
    682 = @(12547, 412)
    683 = @(24257, 689)
    684 = @(57047, 1202)
    685 = @(34441, 1028)
    686 = @(34390, 1139)
    687 = @(31381, 1130)
    688 = @(31925, 1345)
    689 = @(15701, 1061)
    690 = @(27472, 1524)
    691 = @(61944, 2795)
    692 = @(41018, 2646)
    693 = @(38253, 2652)
    694 = @(36096, 2710)
}
foreach ($row in $fgUpdates.Keys) {
    $pair = $fgUpdates[$row]
    $ws.Cells.Item($row, 6).Value = $pair[0]
    $ws.Cells.Item($row, 7).Value = $pair[1]
}

# Append new daily rows (A:G)
$newRows = @{
    695 = @(44589, 987063, 33668, 15118, 17796, 30570, 2620)
    696 = @(44590, 1000345, 26643, 13282, 17817, 13604, 1683)
    697 = @(44591, 1008141, 15798, 7796, 17830, 17449, 1995)
}
foreach ($row in $newRows.Keys) {
    $vals = $newRows[$row]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 1).NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
}
